# Generate Report for Handoff
# Update status text from "In Translation" to "Ready for handoff" and
# refresh the "Latest ... Datetime" timestamps to reflect the new handoff
# report generation time, widening the status columns to fit the new text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# --- Timestamps ---
$overview.Range("G2").Value = "2016-08-21 12:45:39"
$dede.Range("H2").Value = "2016-08-21 12:45:39"
$zhcn.Range("H2").Value = "2016-08-21 12:45:35"

# --- Widen the status columns to fit "Ready for handoff" ---
# (ColumnWidth is quantized to the host's pixel grid, same as real Excel;
# 16.3 is the nearest input that lands on the target rendered width.)
$overview.Range("E1").ColumnWidth = 16.3
$overview.Range("F1").ColumnWidth = 16.3
$zhcn.Range("C1").ColumnWidth = 16.3
$dede.Range("C1").ColumnWidth = 16.3
